$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36-51 down to 37-52.
$ws.Rows("36:36").Insert()

# Populate the newly inserted row 36 with the new weekly entry.
# Most descriptive columns mirror the row that used to be at 36 (now at 37),
# only the measurement columns (D, M, N, O, P, Q, S) carry new data.
$ws.Cells.Item(36, 1).Value = 4
$ws.Cells.Item(36, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(36, 3).Value = "Los Lagos"
$ws.Cells.Item(36, 4).Value = 45072
$ws.Cells.Item(36, 5).Value = 10
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100104
$ws.Cells.Item(36, 8).Value = "Frutos de pepita"
$ws.Cells.Item(36, 9).Value = 100104001
$ws.Cells.Item(36, 10).Value = "Granada"
$ws.Cells.Item(36, 11).Value = "Wonderfull"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 200
$ws.Cells.Item(36, 14).Value = 16000
$ws.Cells.Item(36, 15).Value = 17000
$ws.Cells.Item(36, 16).Value = 16500
$ws.Cells.Item(36, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(36, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(36, 19).Value = 1100
$ws.Cells.Item(36, 20).Value = 15
